$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2180685358255452
$ws.Range("C2").Value = 0.5545171339563862
$ws.Range("J2").Value = 0.009345794392523364
$ws.Range("P2").Value = 0.1433021806853582
$ws.Range("S2").Value = 0.07476635514018691
$ws.Range("B3").Value = 0.005208333333333333
$ws.Range("C3").Value = 0.05208333333333334
$ws.Range("J3").Value = 0.01041666666666667
$ws.Range("P3").Value = 0.7395833333333334
$ws.Range("S3").Value = 0.1927083333333333
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.7560975609756098
$ws.Range("S4").Value = 0.1707317073170732
$ws.Range("B6").Value = 0.08417508417508418
$ws.Range("F6").Value = 0.06734006734006734
$ws.Range("J6").Value = 0.2996632996632997
$ws.Range("O6").Value = 0.02356902356902357
$ws.Range("Q6").Value = 0.1279461279461279
$ws.Range("R6").Value = 0.04713804713804714
$ws.Range("S6").Value = 0.3501683501683502
$ws.Range("B7").Value = 0.1037344398340249
$ws.Range("D7").Value = 0.02074688796680498
$ws.Range("F7").Value = 0.05809128630705394
$ws.Range("J7").Value = 0.1244813278008299
$ws.Range("O7").Value = 0.01244813278008299
$ws.Range("Q7").Value = 0.1742738589211618
$ws.Range("R7").Value = 0.07053941908713693
$ws.Range("S7").Value = 0.4356846473029046
$ws.Range("B8").Value = 0.09090909090909091
$ws.Range("D8").Value = 0.01353965183752418
$ws.Range("E8").Value = 0.001934235976789168
$ws.Range("F8").Value = 0.05415860735009671
$ws.Range("J8").Value = 0.1005802707930367
$ws.Range("O8").Value = 0.02127659574468085
$ws.Range("Q8").Value = 0.1702127659574468
$ws.Range("R8").Value = 0.1353965183752418
$ws.Range("S8").Value = 0.4119922630560928
$ws.Range("B9").Value = 0.08764940239043825
$ws.Range("D9").Value = 0.0199203187250996
$ws.Range("F9").Value = 0.07171314741035857
$ws.Range("J9").Value = 0.151394422310757
$ws.Range("O9").Value = 0.02390438247011952
$ws.Range("Q9").Value = 0.1593625498007968
$ws.Range("R9").Value = 0.08366533864541832
$ws.Range("S9").Value = 0.4023904382470119
$ws.Range("B10").Value = 0.09793420045906656
$ws.Range("D10").Value = 0.01912777352716144
$ws.Range("E10").Value = 0.0007651109410864575
$ws.Range("F10").Value = 0.0864575363427697
$ws.Range("J10").Value = 0.1048201989288447
$ws.Range("O10").Value = 0.02142310635042081
$ws.Range("Q10").Value = 0.2058148431522571
$ws.Range("R10").Value = 0.09487375669472073
$ws.Range("S10").Value = 0.3687834736036725
$ws.Range("G11").Value = 0.1150442477876106
$ws.Range("J11").Value = 0.07079646017699115
$ws.Range("K11").Value = 0.1504424778761062
$ws.Range("L11").Value = 0.6460176991150443
$ws.Range("S11").Value = 0.01769911504424779
$ws.Range("G12").Value = 0.771551724137931
$ws.Range("J12").Value = 0.1681034482758621
$ws.Range("K12").Value = 0.004310344827586207
$ws.Range("L12").Value = 0.03879310344827586
$ws.Range("S12").Value = 0.01724137931034483
$ws.Range("G13").Value = 0.6538461538461539
$ws.Range("J13").Value = 0.2692307692307692
$ws.Range("S13").Value = 0.07692307692307693
$ws.Range("F15").Value = 0.0409556313993174
$ws.Range("H15").Value = 0.174061433447099
$ws.Range("I15").Value = 0.05460750853242321
$ws.Range("J15").Value = 0.3242320819112628
$ws.Range("K15").Value = 0.09215017064846416
$ws.Range("M15").Value = 0.01023890784982935
$ws.Range("O15").Value = 0.09215017064846416
$ws.Range("S15").Value = 0.2116040955631399
$ws.Range("F16").Value = 0.02857142857142857
$ws.Range("H16").Value = 0.1761904761904762
$ws.Range("J16").Value = 0.3619047619047619
$ws.Range("K16").Value = 0.1095238095238095
$ws.Range("M16").Value = 0.004761904761904762
$ws.Range("O16").Value = 0.0761904761904762
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.03617021276595744
$ws.Range("H17").Value = 0.1872340425531915
$ws.Range("I17").Value = 0.1127659574468085
$ws.Range("J17").Value = 0.3808510638297872
$ws.Range("K17").Value = 0.1
$ws.Range("M17").Value = 0.02978723404255319
$ws.Range("O17").Value = 0.07872340425531915
$ws.Range("S17").Value = 0.07446808510638298
$ws.Range("F18").Value = 0.02419354838709677
$ws.Range("H18").Value = 0.1774193548387097
$ws.Range("I18").Value = 0.1088709677419355
$ws.Range("J18").Value = 0.3790322580645161
$ws.Range("K18").Value = 0.1169354838709677
$ws.Range("M18").Value = 0.02016129032258064
$ws.Range("O18").Value = 0.08870967741935484
$ws.Range("S18").Value = 0.0846774193548387
$ws.Range("F19").Value = 0.02814814814814815
$ws.Range("H19").Value = 0.2207407407407407
$ws.Range("I19").Value = 0.1007407407407407
$ws.Range("J19").Value = 0.3340740740740741
$ws.Range("K19").Value = 0.1177777777777778
$ws.Range("M19").Value = 0.0237037037037037
$ws.Range("N19").Value = 0.001481481481481481
$ws.Range("O19").Value = 0.07629629629629629
$ws.Range("S19").Value = 0.09703703703703703
